# Season record columns (Wins / Losses / Ties) appended after the existing
# "Unnamed: 28" column (AC), for the roster table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers, styled like the rest of the
#     header row (bold font, thin box border, centered/top aligned). ---
$headerRange = $ws.Range("AD1:AF1")

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous (thin box border)

# --- Data rows (2..48): every player row gets the team's season record. ---
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 67   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
